# Additional time worked by students to appear for test.
# Update the TIME (column D) entries for the affected days; TOTAL (column F,
# F = D*E) and the grand-total row 162 are formulas, so they recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value  = 2
$ws.Range("D147").Value = 4
$ws.Range("D148").Value = 2
$ws.Range("D149").Value = 3
$ws.Range("D150").Value = 3
$ws.Range("D151").Value = 2

# Restore the sheet's view state: gridlines remain visible, the window is
# scrolled back towards the top of the sheet and the selection moves to E3.
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 130
$win.ScrollColumn = 1
$ws.Range("E3").Select() | Out-Null
